# Weekly update: a new daily price record is inserted at row 10 (pushing the
# existing historical rows down by one), for
# "Hortaliza, Vega Central Mapocho de Santiago - Achicoria".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10. All rows 10..61 shift down to
# 11..62, and the used range grows from A1:R61 to A1:R62.
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with the latest weekly record. The
# non-date/volume fields mirror the same market/category/quality metadata as
# the surrounding rows.
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 45069
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = 100112010
$ws.Cells.Item(10, 7).Value = "Achicoria"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 70
$ws.Cells.Item(10, 11).Value = 7000
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 13).Value = 7000
$ws.Cells.Item(10, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(10, 16).Value = 438
$ws.Cells.Item(10, 17).Value = 16
$ws.Cells.Item(10, 18).Value = "Hortaliza"
